$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 69 ("「外務省」" entry) entirely; all rows below shift up by one.
$ws.Rows.Item(69).Delete()
